{"js": "// The author placed the cursor right after \"...doesn't work und\" (inside the\n// paragraph that begins \"Goes back and searches the previous neighbours...\")\n// and kept typing / editing there -- which is exactly where Word's automatic\n// \"_GoBack\" bookmark (last edit position) ends up. That causes the run that\n// used to read:\n//   \". Furthermore, this approach doesn't work under the current framework\n//    since we are only ever using a single source node at each level. Let's\n//    come back to this later.\"\n// to be split in two, with the _GoBack bookmark sitting between the pieces.\n// The author also added a new, empty trailing paragraph at the very end of\n// the document (e.g. by pressing Enter after the last paragraph).\n\nconst body = context.document.body;\n\n// Word only ever keeps a single \"_GoBack\" bookmark -- remove the existing one\n// before re-inserting it at the new (split) location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Find the split point inside the run's text.\nconst results = body.search(\"this approach doesn't work und\", { matchCase: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the text to split\");\n}\n\nconst splitPoint = results.items[0].getRange(\"After\");\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Add a new empty paragraph at the very end of the document body.\nbody.insertParagraph(\"\", \"End\");\nawait context.sync();\n", "ps1": "# The author placed the cursor right after \"...doesn't work und\" (inside the\n# paragraph that begins \"Goes back and searches the previous neighbours...\")\n# and kept typing / editing there -- which is exactly where Word's automatic\n# \"_GoBack\" bookmark (last edit position) ends up. That causes the run that\n# used to read:\n#   \". Furthermore, this approach doesn't work under the current framework\n#    since we are only ever using a single source node at each level. Let's\n#    come back to this later.\"\n# to be split in two, with the _GoBack bookmark sitting between the pieces.\n# The author also added a new, empty trailing paragraph at the very end of\n# the document (e.g. by pressing Enter after the last paragraph).\n\n$d = $word.ActiveDocument\n\n# Word only ever keeps a single \"_GoBack\" bookmark -- remove the existing one\n# before re-inserting it at the new (split) location.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Find the split point inside the run's text.\n$found = $d.Content\n$ok = $found.Find.Execute(\"this approach doesn't work und\")\n\n# Insert the (now moved) bookmark right after the located text.\n$bmRange = $d.Range($found.End, $found.End)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n# Add a new empty paragraph at the very end of the document.\n$d.Content.InsertParagraphAfter()\n"}
